# form6_report.xlsx aesthetics + backend-integration pass
# - renames the sheet
# - adds a "Проверяющий" signature line and a trailing section
#   (проверку проводил / сноска) below the checklist
# - styles the whole checklist: bold titles, grey header band,
#   yellow "итоговая оценка" band, thin borders around the table
# - widens the columns to fit the Russian text

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------
# 1. Sheet name
# ---------------------------------------------------------------
$ws.Name = "Отчет проверки"

# ---------------------------------------------------------------
# 2. Column widths (character widths; Excel adds ~0.8333 internally
#    to get the stored xlsx width, so we compensate here)
# ---------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 7.1666666666666666   # -> 8
$ws.Columns.Item(2).ColumnWidth = 59.1666666666666666  # -> 60
$ws.Columns.Item(3).ColumnWidth = 11.1666666666666666  # -> 12
$ws.Columns.Item(4).ColumnWidth = 11.1666666666666666  # -> 12
$ws.Columns.Item(5).ColumnWidth = 39.1666666666666666  # -> 40

# ---------------------------------------------------------------
# 3. New / changed cell values
# ---------------------------------------------------------------
$ws.Range("C2").Value = "Дата проведения проверки:"
$ws.Range("A3").Value = "Проверяющий: __________________ __________________ ________________________"

# the old "Итоговая оценка" score moves from C17 to E17
$ws.Range("C17").Value = ""
$ws.Range("E17").Value = 0

$ws.Range("A20").Value = "Проверку проводил _____________________ ______________ _______________________"
$ws.Range("A21").Value = "должность подпись расшифровка подпись"
$ws.Range("A24").Value = "* Итоговая оценка структурному подразделению проставляется проверяющим при выявлении одного и того же несоответствия 2 раза в размере «3 балла», при выявлении одного и того же несоответствия более 2 раз - «2 балла». При отсутствии повторяющихся несоответствий в ходе проведения проверки данная графа проверяющим не заполняется."

# ---------------------------------------------------------------
# 4. Styling
# ---------------------------------------------------------------

# "ЭМО" label - bold, size 14
$r = $ws.Range("A2")
$r.Font.Bold = $true
$r.Font.Size = 14

# "Дата проведения проверки:" / "Проверяющий: ..." - bold, size 12
$r = $ws.Range("C2")
$r.Font.Bold = $true
$r.Font.Size = 12
$r = $ws.Range("A3")
$r.Font.Bold = $true
$r.Font.Size = 12

# Title row - bold size 14, bordered, centered, merged across A5:E5
$r = $ws.Range("A5:E5")
$r.Font.Bold = $true
$r.Font.Size = 14
$r.Borders.LineStyle = 1
$r.HorizontalAlignment = -4108
$r.VerticalAlignment = -4108
$r.WrapText = $true
$r.Merge()

# Header band (A6:E7) - bold, grey fill, bordered, centered + wrapped
$r = $ws.Range("A6:E7")
$r.Font.Bold = $true
$r.Interior.Color = 13882323
$r.Borders.LineStyle = 1
$r.HorizontalAlignment = -4108
$r.VerticalAlignment = -4108
$r.WrapText = $true

# B6 / E6 left-aligned instead of centered
$r = $ws.Range("B6")
$r.HorizontalAlignment = -4131
$r = $ws.Range("E6")
$r.HorizontalAlignment = -4131

# Data table - column A (№, letters) + C:D (checkboxes) centered/bordered
$r = $ws.Range("A8:A14")
$r.Borders.LineStyle = 1
$r.HorizontalAlignment = -4108
$r.VerticalAlignment = -4108
$r.WrapText = $true

$r = $ws.Range("C9:D14")
$r.Borders.LineStyle = 1
$r.HorizontalAlignment = -4108
$r.VerticalAlignment = -4108
$r.WrapText = $true

$r = $ws.Range("C15")
$r.Borders.LineStyle = 1
$r.HorizontalAlignment = -4108
$r.VerticalAlignment = -4108
$r.WrapText = $true

# Data table - column B (criteria text) + E (comments) left-aligned/bordered
$r = $ws.Range("B8:B15")
$r.Borders.LineStyle = 1
$r.HorizontalAlignment = -4131
$r.VerticalAlignment = -4108
$r.WrapText = $true

$r = $ws.Range("E9:E14")
$r.Borders.LineStyle = 1
$r.HorizontalAlignment = -4131
$r.VerticalAlignment = -4108
$r.WrapText = $true

$r = $ws.Range("E17")
$r.Borders.LineStyle = 1
$r.HorizontalAlignment = -4131
$r.VerticalAlignment = -4108
$r.WrapText = $true

# "Итоговая оценка структурному подразделению" - bold 12, yellow fill, merged B17:D17
$r = $ws.Range("B17:D17")
$r.Font.Bold = $true
$r.Font.Size = 12
$r.Interior.Color = 65535
$r.Borders.LineStyle = 1
$r.HorizontalAlignment = -4131
$r.VerticalAlignment = -4108
$r.WrapText = $true
$r.Merge()

# Signature lines (A20, A21) - bordered, centered, wrapped, tall rows, merged
$r = $ws.Range("A20:E20")
$r.Borders.LineStyle = 1
$r.HorizontalAlignment = -4108
$r.VerticalAlignment = -4108
$r.WrapText = $true
$ws.Rows.Item(20).RowHeight = 25
$r.Merge()

$r = $ws.Range("A21:E21")
$r.Borders.LineStyle = 1
$r.HorizontalAlignment = -4108
$r.VerticalAlignment = -4108
$r.WrapText = $true
$ws.Rows.Item(21).RowHeight = 25
$r.Merge()

# filler cells under the signature lines, left aligned, no wrap, no border
$r = $ws.Range("B20:E20")
$r.HorizontalAlignment = -4131
$r.VerticalAlignment = -4108
$r = $ws.Range("B21:E21")
$r.HorizontalAlignment = -4131
$r.VerticalAlignment = -4108

# Footnote - bold 12, yellow fill, bordered, centered, wrapped, tall row, merged
$r = $ws.Range("A24:E26")
$r.Font.Bold = $true
$r.Font.Size = 12
$r.Interior.Color = 65535
$r.Borders.LineStyle = 1
$r.HorizontalAlignment = -4108
$r.VerticalAlignment = -4108
$r.WrapText = $true
$ws.Rows.Item(24).RowHeight = 60
$r.Merge()
